# Append a new batch of performance-test rows (IModel = the new Wraith_MultiMulti
# path, E = 12366) to the "Performance Results" sheet, followed by the usual
# "Min Time" roll-up row, mirroring the existing layout used for every other
# batch already in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$imodelPath = "d:\js\s\imodeljs-core/test-apps/testbed/frontend/performance/imodels/Wraith_MultiMulti.ibim"

# Columns: A=IModel B=View C=Flags D=TileLoadingTime E=Scene F=GarbageExecute
#          G=InitCommands H=BackgroundDraw I=SetClips J=OpaqueDraw K=TranslucentDraw
#          L=HiliteDraw M=CompositeDraw N=OverlayDraw O=RenderFrameTime P=glFinish
#          Q=TotalTime R=(unused header dup, mirrors TotalTime)
# Data rows: A="" (empty string), B=IModel path, C="unknown", D="V0", E=12366,
# then F..R are the per-run numeric timings below.
$data = @(
    @(4, 1, 1, 0, 0, 0, 0, 5, 0, 0, 11, 0, 11),
    @(4, 0, 2, 0, 0, 0, 0, 8, 0, 0, 14, 0, 14),
    @(3, 0, 1, 0, 0, 0, 0, 5, 0, 0, 9, 0, 9),
    @(3, 0, 2, 0, 0, 0, 0, 5, 0, 0, 10, 0, 10),
    @(5, 0, 1, 0, 0, 0, 0, 12, 0, 0, 18, 0, 18),
    @(4, 0, 1, 0, 0, 0, 0, 4, 0, 0, 9, 0, 9),
    @(4, 0, 0, 0, 0, 0, 0, 7, 0, 0, 11, 0, 11),
    @(3, 0, 0, 0, 0, 0, 0, 6, 0, 0, 9, 0, 9),
    @(3, 0, 0, 1, 0, 0, 0, 3, 0, 0, 7, 0, 7),
    @(3, 0, 0, 0, 0, 0, 0, 5, 0, 0, 8, 0, 8),
    @(3, 0, 1, 0, 0, 0, 0, 5, 0, 0, 9, 0, 9),
    @(3, 0, 1, 0, 0, 0, 0, 3, 0, 0, 7, 0, 7),
    @(3, 0, 0, 1, 0, 0, 0, 3, 0, 0, 7, 0, 7),
    @(3, 0, 1, 0, 0, 0, 0, 4, 0, 0, 8, 0, 8),
    @(5, 0, 0, 1, 0, 0, 0, 3, 0, 0, 9, 0, 9),
    @(3, 0, 0, 0, 0, 0, 0, 4, 0, 0, 7, 0, 7),
    @(4, 0, 0, 0, 0, 0, 0, 5, 0, 0, 9, 0, 9),
    @(3, 0, 0, 0, 0, 0, 0, 4, 0, 0, 7, 0, 7),
    @(8, 0, 1, 0, 0, 0, 0, 6, 0, 0, 15, 0, 15),
    @(3, 0, 0, 0, 0, 0, 0, 4, 0, 0, 7, 0, 7)
)

$startRow = 487
$row = $startRow
foreach ($values in $data) {
    # A leading apostrophe forces Excel to store an explicit empty TEXT value
    # (a bare "" would just clear the cell instead), matching the other rows
    # in this sheet whose IModel column is blank text.
    $ws.Cells.Item($row, 1).Value = "'"
    $ws.Cells.Item($row, 2).Value = $imodelPath
    $ws.Cells.Item($row, 3).Value = "unknown"
    $ws.Cells.Item($row, 4).Value = "V0"
    $ws.Cells.Item($row, 5).Value = 12366
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($row, 6 + $i).Value = $values[$i]
    }
    $row++
}

# Trailing "Min Time" summary row for this batch (only columns D..R populated,
# matching the other batch roll-up rows already on the sheet).
$ws.Cells.Item($row, 4).Value = "Min Time"
$ws.Cells.Item($row, 5).Value = 12366
$minValues = @(3, 0, 0, 0, 0, 0, 0, 4, 0, 0, 7, 0, 7)
for ($i = 0; $i -lt $minValues.Count; $i++) {
    $ws.Cells.Item($row, 6 + $i).Value = $minValues[$i]
}
